$wb = $excel.ActiveWorkbook

# Helper to write a text value into a cell while preventing Excel's automatic
# locale-based type inference (e.g. turning "02/12/2024" into a date serial
# number). We force the cell to Text format first, assign the value, then
# reset the cell style back to "Normal" so no stray number-format style is
# left behind on the cell (matching the source data which has no style).
function Set-TextValue {
    param($ws, [string]$addr, $value)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet "SPN": append rows 108-114
# ---------------------------------------------------------------------
$wsSPN = $wb.Worksheets.Item("SPN")

$spnRows = @(
    @(108, "SPN", "Arthur Hassuma", 9, "02/12/2024", "06/12/2024", 314457, "11/2024", "02/12/2024", "Resolvido", "Willian Jones"),
    @(109, "SPN", "Arthur Hassuma", 9, "02/12/2024", "06/12/2024", 315282, "11/2024", "02/12/2024", "Resolvido", "Willian Jones"),
    @(110, "SPN", "Arthur Hassuma", 9, "02/12/2024", "06/12/2024", 315511, "11/2024", "02/12/2024", "Resolvido", "Willian Jones"),
    @(111, "SPN", "Arthur Hassuma", 9, "02/12/2024", "06/12/2024", 315663, "11/2024", "02/12/2024", "Resolvido", "Willian Jones"),
    @(112, "SPN", "Higor Cruz",     9, "02/12/2024", "06/12/2024", 315374, "11/2024", "02/12/2024", "Resolvido", "Willian Jones"),
    @(113, "SPN", "Luan Pierry",    9, "02/12/2024", "06/12/2024", 315638, "11/2024", "02/12/2024", "Resolvido", "Willian Jones"),
    @(114, "SPN", "Luan Pierry",    9, "02/12/2024", "06/12/2024", 315817, "11/2024", "02/12/2024", "Resolvido", "Willian Jones")
)

foreach ($row in $spnRows) {
    $r   = $row[0]
    $wsSPN.Range("A$r").Value = $row[1]
    $wsSPN.Range("B$r").Value = $row[2]
    $wsSPN.Range("C$r").Value = $row[3]
    Set-TextValue $wsSPN "D$r" $row[4]
    Set-TextValue $wsSPN "E$r" $row[5]
    $wsSPN.Range("F$r").Value = $row[6]
    Set-TextValue $wsSPN "G$r" $row[7]
    Set-TextValue $wsSPN "H$r" $row[8]
    $wsSPN.Range("I$r").Value = $row[9]
    $wsSPN.Range("J$r").Value = $row[10]
}

# ---------------------------------------------------------------------
# Sheet "ITI": update status of existing rows 123 and 129, then append
# rows 135-158
# ---------------------------------------------------------------------
$wsITI = $wb.Worksheets.Item("ITI")

$wsITI.Range("I123").Value = "Resolvido"
$wsITI.Range("I129").Value = "Resolvido"

$itiRows = @(
    @(135, "ITI", "Erick Silva",       9, "02/12/2024", "06/12/2024", 315595, "11/2024", "02/12/2024", "Resolvido", "Emerson Simette"),
    @(136, "ITI", "Erick Silva",       9, "02/12/2024", "06/12/2024", 315683, "11/2024", "02/12/2024", "Resolvido", "Emerson Simette"),
    @(137, "ITI", "Erick Silva",       9, "02/12/2024", "06/12/2024", 315754, "11/2024", "02/12/2024", "Resolvido", "Emerson Simette"),
    @(138, "ITI", "Gustavo Linpiski",  9, "02/12/2024", "06/12/2024", 315377, "11/2024", "02/12/2024", "Resolvido", "Emerson Simette"),
    @(139, "ITI", "Gustavo Linpiski",  9, "02/12/2024", "06/12/2024", 315966, "11/2024", "02/12/2024", "Resolvido", "Emerson Simette"),
    @(140, "ITI", "Jorgenaldo Reis",   9, "02/12/2024", "06/12/2024", 315952, "11/2024", "02/12/2024", "Pendente",  "Emerson Simette"),
    @(141, "ITI", "Jorgenaldo Reis",   9, "02/12/2024", "06/12/2024", 316151, "11/2024", "02/12/2024", "Pendente",  "Emerson Simette"),
    @(142, "ITI", "Jorgenaldo Reis",   9, "02/12/2024", "06/12/2024", 315807, "11/2024", "02/12/2024", "Resolvido", "Emerson Simette"),
    @(143, "ITI", "Jorgenaldo Reis",   9, "02/12/2024", "06/12/2024", 316110, "11/2024", "02/12/2024", "Resolvido", "Emerson Simette"),
    @(144, "ITI", "Jose Acevedo",      9, "02/12/2024", "06/12/2024", 315183, "11/2024", "02/12/2024", "Pendente",  "Emerson Simette"),
    @(145, "ITI", "Alana Neris",      10, "09/12/2024", "13/12/2024", 316702, "12/2024", "09/12/2024", "Pendente",  "Emerson Simette"),
    @(146, "ITI", "Alana Neris",      10, "09/12/2024", "13/12/2024", 316765, "12/2024", "09/12/2024", "Pendente",  "Emerson Simette"),
    @(147, "ITI", "Alana Neris",      10, "09/12/2024", "13/12/2024", 316993, "12/2024", "09/12/2024", "Pendente",  "Emerson Simette"),
    @(148, "ITI", "Edson Campos",     10, "09/12/2024", "13/12/2024", 315916, "11/2024", "09/12/2024", "Pendente",  "Emerson Simette"),
    @(149, "ITI", "Edson Campos",     10, "09/12/2024", "13/12/2024", 315812, "11/2024", "09/12/2024", "Pendente",  "Emerson Simette"),
    @(150, "ITI", "Erick Silva",      10, "09/12/2024", "13/12/2024", 316732, "12/2024", "09/12/2024", "Pendente",  "Emerson Simette"),
    @(151, "ITI", "Erick Silva",      10, "09/12/2024", "13/12/2024", 316626, "12/2024", "09/12/2024", "Pendente",  "Emerson Simette"),
    @(152, "ITI", "Erick Silva",      10, "09/12/2024", "13/12/2024", 316501, "12/2024", "09/12/2024", "Pendente",  "Emerson Simette"),
    @(153, "ITI", "Gabriel Lopez",    10, "09/12/2024", "13/12/2024", 316114, "11/2024", "09/12/2024", "Pendente",  "Emerson Simette"),
    @(154, "ITI", "Jacyr Popenda",    10, "09/12/2024", "13/12/2024", 316940, "12/2024", "09/12/2024", "Pendente",  "Emerson Simette"),
    @(155, "ITI", "Jorgenaldo Reis",  10, "09/12/2024", "13/12/2024", 315817, "11/2024", "09/12/2024", "Pendente",  "Emerson Simette"),
    @(156, "ITI", "Jose Acevedo",     10, "09/12/2024", "13/12/2024", 316763, "12/2024", "09/12/2024", "Pendente",  "Emerson Simette"),
    @(157, "ITI", "Lourival Moizés",  10, "09/12/2024", "13/12/2024", 315310, "11/2024", "09/12/2024", "Pendente",  "Emerson Simette"),
    @(158, "ITI", "Lourival Moizés",  10, "09/12/2024", "13/12/2024", 315758, "11/2024", "09/12/2024", "Pendente",  "Emerson Simette")
)

foreach ($row in $itiRows) {
    $r   = $row[0]
    $wsITI.Range("A$r").Value = $row[1]
    $wsITI.Range("B$r").Value = $row[2]
    $wsITI.Range("C$r").Value = $row[3]
    Set-TextValue $wsITI "D$r" $row[4]
    Set-TextValue $wsITI "E$r" $row[5]
    $wsITI.Range("F$r").Value = $row[6]
    Set-TextValue $wsITI "G$r" $row[7]
    Set-TextValue $wsITI "H$r" $row[8]
    $wsITI.Range("I$r").Value = $row[9]
    $wsITI.Range("J$r").Value = $row[10]
}
